# Weekly update: shift the Coliflor / Vega Modelo de Temuco price rows
# down by one row and insert the new weeks observation at row 386.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 409 is brand new -- seed the columns that stay constant for every
# "Vega Modelo de Temuco" / "Coliflor" observation before filling the
# shifted values below.
$ws.Range("A409").Value = 10
$ws.Range("B409").Value = "Vega Modelo de Temuco"
$ws.Range("C409").Value = "La Araucanía"
$ws.Range("E409").Value = 9
$ws.Range("F409").Value = 100112008
$ws.Range("G409").Value = "Coliflor"
$ws.Range("H409").Value = "Sin especificar"
$ws.Range("I409").Value = "Primera"
$ws.Range("N409").Value = "`$/unidad"
$ws.Range("Q409").Value = 1
$ws.Range("R409").Value = "Hortaliza"

# Shift the date / price / origin columns down by one row (386->409),
# writing the new top observation into row 386.

$ws.Range("D386").Value = 44746
$ws.Range("J386").Value = 1450
$ws.Range("K386").Value = 1000
$ws.Range("L386").Value = 1200
$ws.Range("M386").Value = 1090
$ws.Range("O386").Value = "Provincia de Cautín"
$ws.Range("P386").Value = 1090

$ws.Range("D387").Value = 44386
$ws.Range("J387").Value = 1500
$ws.Range("K387").Value = 800
$ws.Range("L387").Value = 800
$ws.Range("M387").Value = 800
$ws.Range("O387").Value = "Región Metropolitana"
$ws.Range("P387").Value = 800

$ws.Range("D388").Value = 44690
$ws.Range("J388").Value = 1250
$ws.Range("K388").Value = 1400
$ws.Range("L388").Value = 1400
$ws.Range("M388").Value = 1400
$ws.Range("O388").Value = "Región de O'Higgins"
$ws.Range("P388").Value = 1400

$ws.Range("D389").Value = 44690
$ws.Range("J389").Value = 1550
$ws.Range("K389").Value = 1200
$ws.Range("L389").Value = 1200
$ws.Range("M389").Value = 1200
$ws.Range("O389").Value = "Región del Maule"
$ws.Range("P389").Value = 1200

$ws.Range("D390").Value = 44631
$ws.Range("J390").Value = 1250
$ws.Range("K390").Value = 1200
$ws.Range("L390").Value = 1300
$ws.Range("M390").Value = 1252
$ws.Range("O390").Value = "Región de O'Higgins"
$ws.Range("P390").Value = 1252

$ws.Range("D391").Value = 44582
$ws.Range("J391").Value = 500
$ws.Range("K391").Value = 1000
$ws.Range("L391").Value = 1000
$ws.Range("M391").Value = 1000
$ws.Range("O391").Value = "Provincia de Cautín"
$ws.Range("P391").Value = 1000

$ws.Range("D392").Value = 44307
$ws.Range("J392").Value = 850
$ws.Range("K392").Value = 1000
$ws.Range("L392").Value = 1000
$ws.Range("M392").Value = 1000
$ws.Range("O392").Value = "Región Metropolitana"
$ws.Range("P392").Value = 1000

$ws.Range("D393").Value = 44307
$ws.Range("J393").Value = 1250
$ws.Range("K393").Value = 1000
$ws.Range("L393").Value = 1000
$ws.Range("M393").Value = 1000
$ws.Range("O393").Value = "Región del Maule"
$ws.Range("P393").Value = 1000

$ws.Range("D394").Value = 44344
$ws.Range("J394").Value = 2000
$ws.Range("K394").Value = 1000
$ws.Range("L394").Value = 1000
$ws.Range("M394").Value = 1000
$ws.Range("O394").Value = "Región Metropolitana"
$ws.Range("P394").Value = 1000

$ws.Range("D395").Value = 44433
$ws.Range("J395").Value = 1850
$ws.Range("K395").Value = 800
$ws.Range("L395").Value = 800
$ws.Range("M395").Value = 800
$ws.Range("O395").Value = "Región de O'Higgins"
$ws.Range("P395").Value = 800

$ws.Range("D396").Value = 44707
$ws.Range("J396").Value = 800
$ws.Range("K396").Value = 1200
$ws.Range("L396").Value = 1200
$ws.Range("M396").Value = 1200
$ws.Range("O396").Value = "Provincia de Cautín"
$ws.Range("P396").Value = 1200

$ws.Range("D397").Value = 44707
$ws.Range("J397").Value = 5000
$ws.Range("K397").Value = 1200
$ws.Range("L397").Value = 1200
$ws.Range("M397").Value = 1200
$ws.Range("O397").Value = "Región Metropolitana"
$ws.Range("P397").Value = 1200

$ws.Range("D398").Value = 44707
$ws.Range("J398").Value = 2000
$ws.Range("K398").Value = 1200
$ws.Range("L398").Value = 1200
$ws.Range("M398").Value = 1200
$ws.Range("O398").Value = "Región del Maule"
$ws.Range("P398").Value = 1200

$ws.Range("D399").Value = 44421
$ws.Range("J399").Value = 4300
$ws.Range("K399").Value = 800
$ws.Range("L399").Value = 1000
$ws.Range("M399").Value = 884
$ws.Range("O399").Value = "Región Metropolitana"
$ws.Range("P399").Value = 884

$ws.Range("D400").Value = 44421
$ws.Range("J400").Value = 1450
$ws.Range("K400").Value = 800
$ws.Range("L400").Value = 900
$ws.Range("M400").Value = 845
$ws.Range("O400").Value = "Región de O'Higgins"
$ws.Range("P400").Value = 845

$ws.Range("D401").Value = 44637
$ws.Range("J401").Value = 600
$ws.Range("K401").Value = 1400
$ws.Range("L401").Value = 1400
$ws.Range("M401").Value = 1400
$ws.Range("O401").Value = "Provincia de Cautín"
$ws.Range("P401").Value = 1400

$ws.Range("D402").Value = 44637
$ws.Range("J402").Value = 800
$ws.Range("K402").Value = 1300
$ws.Range("L402").Value = 1300
$ws.Range("M402").Value = 1300
$ws.Range("O402").Value = "Región Metropolitana"
$ws.Range("P402").Value = 1300

$ws.Range("D403").Value = 44195
$ws.Range("J403").Value = 650
$ws.Range("K403").Value = 1100
$ws.Range("L403").Value = 1100
$ws.Range("M403").Value = 1100
$ws.Range("O403").Value = "Región del Maule"
$ws.Range("P403").Value = 1100

$ws.Range("D404").Value = 44442
$ws.Range("J404").Value = 400
$ws.Range("K404").Value = 800
$ws.Range("L404").Value = 900
$ws.Range("M404").Value = 850
$ws.Range("O404").Value = "Región Metropolitana"
$ws.Range("P404").Value = 850

$ws.Range("D405").Value = 44483
$ws.Range("J405").Value = 2700
$ws.Range("K405").Value = 800
$ws.Range("L405").Value = 900
$ws.Range("M405").Value = 844
$ws.Range("O405").Value = "Región Metropolitana"
$ws.Range("P405").Value = 844

$ws.Range("D406").Value = 44483
$ws.Range("J406").Value = 1500
$ws.Range("K406").Value = 800
$ws.Range("L406").Value = 900
$ws.Range("M406").Value = 853
$ws.Range("O406").Value = "Región de O'Higgins"
$ws.Range("P406").Value = 853

$ws.Range("D407").Value = 44188
$ws.Range("J407").Value = 1600
$ws.Range("K407").Value = 1000
$ws.Range("L407").Value = 1100
$ws.Range("M407").Value = 1050
$ws.Range("O407").Value = "Región del Maule"
$ws.Range("P407").Value = 1050

$ws.Range("D408").Value = 44519
$ws.Range("J408").Value = 2050
$ws.Range("K408").Value = 800
$ws.Range("L408").Value = 900
$ws.Range("M408").Value = 861
$ws.Range("O408").Value = "Región del Maule"
$ws.Range("P408").Value = 861

$ws.Range("D409").Value = 44194
$ws.Range("J409").Value = 1550
$ws.Range("K409").Value = 1100
$ws.Range("L409").Value = 1100
$ws.Range("M409").Value = 1100
$ws.Range("O409").Value = "Región del Maule"
$ws.Range("P409").Value = 1100

# Make sure the new date cell carries the same date number-format as the
# rest of column D.
$ws.Range("D409").NumberFormat = $ws.Range("D408").NumberFormat

Write-Output "done"
